$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 37
$ws1.Range("H2").Value = 4.32
$ws1.Range("L2").Value = 1.08

# Row 3
$ws1.Range("D3").Value = 29
$ws1.Range("H3").Value = 4.18
$ws1.Range("L3").Value = 1.11

# Row 4
$ws1.Range("D4").Value = 30
$ws1.Range("H4").Value = 3.11
$ws1.Range("L4").Value = 1.14

# Row 5
$ws1.Range("D5").Value = 32
$ws1.Range("H5").Value = 1.97
$ws1.Range("J5").Value = "Normal"
$ws1.Range("L5").Value = 0.97

# Row 6
$ws1.Range("D6").Value = 32
$ws1.Range("H6").Value = 0.97
$ws1.Range("I6").Value = "Low"
$ws1.Range("L6").Value = 1.12

# Row 7
$ws1.Range("D7").Value = 28
$ws1.Range("L7").Value = 0.91

# Row 8
$ws1.Range("D8").Value = 35
$ws1.Range("L8").Value = 1.09

# Row 9
$ws1.Range("D9").Value = 35
$ws1.Range("L9").Value = 1.2

# Row 10
$ws1.Range("D10").Value = 30
$ws1.Range("L10").Value = 0.96

# Row 11
$ws1.Range("D11").Value = 32
$ws1.Range("L11").Value = 1.19

# Row 12
$ws1.Range("D12").Value = 33
$ws1.Range("L12").Value = 1.11

# Row 13
$ws1.Range("D13").Value = 33
$ws1.Range("L13").Value = 0.84

# Row 14
$ws1.Range("D14").Value = 34
$ws1.Range("L14").Value = 1.12

# Row 15
$ws1.Range("D15").Value = 30
$ws1.Range("L15").Value = 0.89

# Row 16
$ws1.Range("D16").Value = 33
$ws1.Range("L16").Value = 1.13

# Row 17
$ws1.Range("D17").Value = 32
$ws1.Range("L17").Value = 0.84

# --- Sheet: Summary ---
# These cells hold numeric-looking values stored as text, so we prefix
# with an apostrophe to force text entry (keeps cell type as string).
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'523"
$ws2.Range("B10").Value = "'260"
$ws2.Range("B11").Value = "'129"
$ws2.Range("B12").Value = "'37"
$ws2.Range("B14").Value = "'29"
